$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2876
$ws1.Range("F7").Value = 1688
$ws1.Range("F8").Value = 1891
$ws1.Range("F10").Value = 287
$ws1.Range("F11").Value = 772
$ws1.Range("F12").Value = 918
$ws1.Range("F13").Value = 177
$ws1.Range("F17").Value = 56
$ws1.Range("F19").Value = 6858
$ws1.Range("F20").Value = 259
$ws1.Range("F21").Value = 1653
$ws1.Range("F23").Value = 183
$ws1.Range("F25").Value = 320
$ws1.Range("F26").Value = 277
$ws1.Range("F29").Value = 921
$ws1.Range("F31").Value = 101
$ws1.Range("F34").Value = 1923
$ws1.Range("F37").Value = 234
$ws1.Range("F38").Value = 27
$ws1.Range("F40").Value = 236
$ws1.Range("F42").Value = 173

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 5

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 2876
$ws4.Range("F10").Value = 1688
$ws4.Range("F11").Value = 1891
$ws4.Range("F13").Value = 287
$ws4.Range("F14").Value = 772
$ws4.Range("F16").Value = 918
$ws4.Range("F17").Value = 177
$ws4.Range("F20").Value = 56
$ws4.Range("F22").Value = 6858
$ws4.Range("F23").Value = 259
$ws4.Range("F24").Value = 1653
$ws4.Range("F25").Value = 5
$ws4.Range("F27").Value = 183
$ws4.Range("F29").Value = 320
$ws4.Range("F30").Value = 277
$ws4.Range("F33").Value = 921
$ws4.Range("F35").Value = 101
$ws4.Range("F38").Value = 1923
$ws4.Range("F41").Value = 234
$ws4.Range("F42").Value = 27
$ws4.Range("F44").Value = 236
$ws4.Range("F49").Value = 173
